$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 6 (second iteration of the gain measurements) ---
$ws.Range("C6").Value = 5275
$ws.Range("D6").Value = 2000
$ws.Range("E6").Formula = "=D6/C6"
$ws.Range("D6").NumberFormat = "#,##0\ ""mm"""

$ws.Range("G6").Value = 5301
$ws.Range("H6").Value = 2000
$ws.Range("I6").Formula = "=H6/G6"
$ws.Range("H6").NumberFormat = "#,##0\ ""mm"""

# --- New "BF cap et avance" rotation-gain block (rows 17-18) ---
$ws.Range("B17").Value = "ancien gain rot"
$ws.Range("C17").Value = "theorique"
$ws.Range("D17").Value = "mesuree"
$ws.Range("E17").Value = "new gain"

$ws.Range("B18").Value = 0.006
$ws.Range("C18").Value = 216.56
$ws.Range("D18").Value = 180
$ws.Range("E18").Formula = "=B18*D18/C18"

# --- Column H should mirror column D's width (same bestFit custom width) ---
$ws.Columns("H").ColumnWidth = $ws.Columns("D").ColumnWidth

# --- Move / record the active selection where the user ended up editing ---
$ws.Range("C19").Select()
